# Update the student id column (A1:A90) from 7001..7090 to plain 1..90
# (commit: "update: csv and xlsx file for the student id to int")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 1; $i -le 90; $i++) {
    $ws.Cells.Item($i, 1).Value = $i
}

# Reflect the saved selection/view state from the diff: the sheet was left
# with A1:A90 selected and scrolled back to the top of the sheet.
$ws.Range("A1:A90").Select()
